$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, [string]$value)
    # Prefix with an apostrophe so Excel always stores the value as text
    # (prevents numeric-looking strings like "296.44" or "39.458.26" from
    # being converted into numeric cells), then reset the style so no
    # quote-prefix formatting/style index leaks into the saved file.
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "39.458.26"
$ws.Range("E2").Value = "  -1.19%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.220.88"
$ws.Range("E3").Value = "  -4.54%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
Set-TextCell "D5" "296.44"
$ws.Range("E5").Value = "  -3.79%  "

# Row 6 - Solana
Set-TextCell "D6" "81.62"
$ws.Range("E6").Value = "  -3.63%  "

# Row 7 - XRP
Set-TextCell "D7" "0.513"
$ws.Range("E7").Value = "  -3.54%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.470"
$ws.Range("E9").Value = "  -2.92%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.0776"
$ws.Range("E10").Value = "  -4.58%  "

# Row 11 - Avalanche
Set-TextCell "D11" "29.78"
$ws.Range("E11").Value = "  -0.68%  "

# Row 12 - OKB
Set-TextCell "D12" "46.84"
$ws.Range("E12").Value = "  -11.15%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell "D14" "2.556.03"
$ws.Range("E14").Value = "  -4.83%  "

# Row 15 - Polkadot
Set-TextCell "D15" "6.27"
$ws.Range("E15").Value = "  -2.05%  "

# Row 16 - Chainlink
Set-TextCell "D16" "14.03"
$ws.Range("E16").Value = "  -4.31%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.211.50"
$ws.Range("E17").Value = "  -4.73%  "

# Row 18 - Polygon
Set-TextCell "D18" "0.716"
$ws.Range("E18").Value = "  -4.98%  "

# Row 19 - WrappedBTC
Set-TextCell "D19" "39.360.96"
$ws.Range("E19").Value = "  -1.43%  "

# Row 20 - ShibaInu
Set-TextCell "D20" "0.0₃0874"
$ws.Range("E20").Value = "  -3.15%  "

# Row 21 - Uniswap
Set-TextCell "D21" "5.74"
$ws.Range("E21").Value = "  -5.34%  "

# Row 22 - Litecoin
Set-TextCell "D22" "64.84"
$ws.Range("E22").Value = "  -3.97%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextCell "D23" "10.30"
$ws.Range("E23").Value = "  -2.89%  "

# Row 24 - BitcoinCash
Set-TextCell "D24" "229.90"
$ws.Range("E24").Value = "  -2.23%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  -0.10%  "

# Row 26 - PancakeSwap
Set-TextCell "D26" "2.41"
$ws.Range("E26").Value = "  -5.56%  "

# Row 27 - ImmutableX
Set-TextCell "D27" "1.81"
$ws.Range("E27").Value = "  +0.64%  "

# Row 28 - EthereumClassic
Set-TextCell "D28" "22.68"
$ws.Range("E28").Value = "  -2.58%  "

# Row 29 - Toncoin
Set-TextCell "D29" "2.18"
$ws.Range("E29").Value = "  -0.99%  "

# Row 30 - Cosmos
Set-TextCell "D30" "9.12"
$ws.Range("E30").Value = "  -1.57%  "

# Row 31 - Monero
Set-TextCell "D31" "149.60"
$ws.Range("E31").Value = "  -1.89%  "

# Row 32 - InjectiveProtocol
Set-TextCell "D32" "31.86"
$ws.Range("E32").Value = "  -9.06%  "

# Row 33 - FirstDigitalUSD
Set-TextCell "D33" "0.999"
$ws.Range("E33").Value = "  -0.28%  "

# Row 34 - Filecoin
Set-TextCell "D34" "4.81"
$ws.Range("E34").Value = "  -5.51%  "

# Row 35 - Hedera
Set-TextCell "D35" "0.0697"
$ws.Range("E35").Value = "  -3.25%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  -4.17%  "

# Row 37 - Stellar
$ws.Range("E37").Value = "  -2.58%  "

# Row 38 - Celestia
Set-TextCell "D38" "15.67"
$ws.Range("E38").Value = "  +0.40%  "

# Row 39 - Kaspa
Set-TextCell "D39" "0.0959"
$ws.Range("E39").Value = "  -3.76%  "

# Row 40 - LidoDAOToken
Set-TextCell "D40" "2.65"
$ws.Range("E40").Value = "  -3.68%  "

# Row 41 - ARBITRUM
Set-TextCell "D41" "1.66"
$ws.Range("E41").Value = "  -2.37%  "

# Row 42 - RenderToken
Set-TextCell "D42" "3.64"
$ws.Range("E42").Value = "  -4.72%  "

# Row 43 - Maker
Set-TextCell "D43" "1.909.58"
$ws.Range("E43").Value = "  -2.03%  "

# Row 44 - ApeXProtocol
Set-TextCell "D44" "2.05"
$ws.Range("E44").Value = "  -9.08%  "

# Row 45 - VeChain
Set-TextCell "D45" "0.0260"
$ws.Range("E45").Value = "  -2.16%  "

# Row 46 - FraxShare
Set-TextCell "D46" "9.14"
$ws.Range("E46").Value = "  -2.39%  "

# Row 47 - EnergySwap
Set-TextCell "D47" "16.35"
$ws.Range("E47").Value = "  -6.81%  "

# Row 48 - NEARProtocol
Set-TextCell "D48" "2.63"
$ws.Range("E48").Value = "  -1.62%  "

# Row 49 - RocketPoolETH
Set-TextCell "D49" "2.423.82"
$ws.Range("E49").Value = "  -5.08%  "

# Row 50 - BitcoinSV
Set-TextCell "D50" "71.15"
$ws.Range("E50").Value = "  +0.98%  "

# Row 51 - Aave
Set-TextCell "D51" "88.15"
$ws.Range("E51").Value = "  -4.89%  "
